$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Colour used by the workbook's existing "file name" hyperlinks
# (FF6495ED, i.e. cornflower blue) expressed as a BGR OLE colour value,
# plus single-underline, so that newly (re)created hyperlinks keep the
# same look as the ones that were already in the sheet.
$hlColor = 15570276   # RGB(100,149,237) -> 0x00ED9564 (BGR)
$hlUnderlineSingle = 2

# ---- Remove every existing hyperlink on each sheet so we can re-add them ----
# (in final order) with fresh, sequential relationship ids that match the
# target layout. (Range.Hyperlinks.Delete() clears every hyperlink on the
# whole sheet in this engine, regardless of which range it is called on.)
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsDeDe.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.md", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md", "", "", "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8283f95e7cbd1d5e4478644934b2e5ff383053d4/.localization-config", "", "", ".localization-config") | Out-Null

$wsOverview.Range("A2:A4").Font.Color = $hlColor
$wsOverview.Range("A2:A4").Font.Underline = $hlUnderlineSingle

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn.Range("A2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-02-17 10:08:47"
$wsZhCn.Range("G2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-02-17 10:08:47"
$wsZhCn.Range("G3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H3").Value = "Include"

$wsZhCn.Range("A4").Value = ".localization-config"
$wsZhCn.Range("B4").Value = "Not to be localized"
$wsZhCn.Range("D4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("G4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H4").Value = "Ignored"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.md", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md", "", "", "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8283f95e7cbd1d5e4478644934b2e5ff383053d4/.localization-config", "", "", ".localization-config") | Out-Null

$wsZhCn.Range("A2").Font.Color = $hlColor
$wsZhCn.Range("A2").Font.Underline = $hlUnderlineSingle
$wsZhCn.Range("A3").Font.Color = $hlColor
$wsZhCn.Range("A3").Font.Underline = $hlUnderlineSingle
$wsZhCn.Range("A4").Font.Color = $hlColor
$wsZhCn.Range("A4").Font.Underline = $hlUnderlineSingle
$wsZhCn.Range("C2").Font.Color = $hlColor
$wsZhCn.Range("C2").Font.Underline = $hlUnderlineSingle
$wsZhCn.Range("C3").Font.Color = $hlColor
$wsZhCn.Range("C3").Font.Underline = $hlUnderlineSingle

$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe.Range("A2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md"
$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-02-17 10:08:58"
$wsDeDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md"
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-02-17 10:08:58"
$wsDeDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H3").Value = "Include"

$wsDeDe.Range("A4").Value = ".localization-config"
$wsDeDe.Range("B4").Value = "Not to be localized"
$wsDeDe.Range("D4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H4").Value = "Ignored"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.md", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md", "", "", "ffffd646bfcc-de81-48ca-9a72-42d51bb74aa6.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTest/oltest/blob/b02f3e3848a83d70c4fb3eb84658976e950a945a/e2e/a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf", "", "", "a9de59e0-3756-49d9-b0b7-ac1152a258f9.bcb779f2f44eea98f9e969697df2388c6fe071a3.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/8283f95e7cbd1d5e4478644934b2e5ff383053d4/.localization-config", "", "", ".localization-config") | Out-Null

$wsDeDe.Range("A2").Font.Color = $hlColor
$wsDeDe.Range("A2").Font.Underline = $hlUnderlineSingle
$wsDeDe.Range("A3").Font.Color = $hlColor
$wsDeDe.Range("A3").Font.Underline = $hlUnderlineSingle
$wsDeDe.Range("A4").Font.Color = $hlColor
$wsDeDe.Range("A4").Font.Underline = $hlUnderlineSingle
$wsDeDe.Range("C2").Font.Color = $hlColor
$wsDeDe.Range("C2").Font.Underline = $hlUnderlineSingle
$wsDeDe.Range("C3").Font.Color = $hlColor
$wsDeDe.Range("C3").Font.Underline = $hlUnderlineSingle

$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
